$d = $word.ActiveDocument

# Change 1: add sentence about split plot design with four replicates
$rng1 = $d.Content
$found1 = $rng1.Find.Execute(
    "findings from the 2015 season. The changes are detailed following.",
    $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)
if ($found1) {
    $rng1.Text = "findings from the 2015 season. Both seasons consisted of split plot design with four replicates where irrigation was the main plot and N rate was the split plot treatent. The changes are detailed following."
}

# Change 2: remove leading sentence from 2015 Dry Season paragraph and append new sentence at the end
$rng2 = $d.Content
$found2 = $rng2.Find.Execute(
    "Split plot design with irrigation as the main plot treatment and N rate as the split plot treatment. The main plot size was 12m x 12m (144 sq m), with a sub-plot size of 5m x 5m (25 sq m). Replication size was 12m x 24m (288 sq m) with a buffer of 1m per sub plot for a whole experiment size of 1,152 sq m.",
    $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)
if ($found2) {
    $rng2.Text = "The main plot size was 12m x 12m (144 sq m), with a sub-plot size of 5m x 5m (25 sq m). Replication size was 12m x 24m (288 sq m) with a buffer of 1m per sub plot for a whole experiment size of 1,152 sq m. The main plot treatments were alternate wetting and drying (AWD) and flooded or farmers' practice."
}

# Change 3: rephrase start of 2016 Dry Season paragraph
$rng3 = $d.Content
$found3 = $rng3.Find.Execute(
    "A split plot design was used again with irrigation as the main plot treatment and N rate as the split treatment. However, the plot size increased and due to these changes,",
    $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)
if ($found3) {
    $rng3.Text = "In 2016 dry season the plot size was increased and due to these changes,"
}

Write-Output "found1=$found1 found2=$found2 found3=$found3"
